$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date serial number that was bumped by one day
# (46061 -> 46062) for every data row, from row 2 through row 550.
$ws.Range("C2:C550").Value = 46062
